$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Remove the 16 rows belonging to the first worker (HENRY ANTONIO
#    GUERRERO ROMERO / CC 73190080). This shifts the second worker's
#    block (EVA YULIANIS GUERRA VARGAS, previously rows 32-50) up to
#    rows 16-34.
# ------------------------------------------------------------------
$ws.Rows("16:31").Delete()

# ------------------------------------------------------------------
# 2) Insert one extra row before the last (bottom-bordered) row of the
#    table so that the table grows from 19 to 20 period rows (adding
#    period 2508). The new row inherits the "middle" row formatting
#    from the row above it, while the old last row (now pushed down)
#    keeps its special bottom-border formatting.
# ------------------------------------------------------------------
$ws.Rows("34:34").Insert()

# ------------------------------------------------------------------
# 3) Fill in the worker/period table, rows 16-35, all for the single
#    remaining worker (EVA YULIANIS GUERRA VARGAS, CC 1007275831),
#    periods 2401..2508 in ascending order.
# ------------------------------------------------------------------
$periods = @("2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412","2501","2502","2503","2504","2505","2506","2507","2508")

for ($i = 0; $i -lt 20; $i++) {
    $r = 16 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "1007275831"
    $ws.Range("D$r").Value = "EVA YULIANIS GUERRA VARGAS"
    $ws.Range("E$r").Value = $periods[$i]
    if ($r -eq 16) {
        $ws.Range("F$r").Value = 6933
    } else {
        $ws.Range("F$r").Value = 52000
    }
    $ws.Range("G$r").Value = 1300000
}

# ------------------------------------------------------------------
# 4) Update the summary header values.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 994933
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 20

# ------------------------------------------------------------------
# 5) Update the signature block (now rows 40-41 after the row
#    deletions/insertion above).
# ------------------------------------------------------------------
$ws.Range("B40").Value = "___________________________________"
$ws.Range("H40").Value = "___________________________________"
$ws.Range("B41").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H41").Value = "FIRMA DEL REPRESENTANTE LEGAL"
